$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 / column C holds the numeric condition value for rule "R30".
# Update it from 18 to 1.
$ws.Range("C10").Value = 1
